$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 589.8889
$ws.Range("I9").Value = 651.125
$ws.Range("K9").Value = 651.125
$ws.Range("M9").Value = -482.125

$ws.Range("H33").Value = 163.45454
$ws.Range("I33").Value = 159.8
$ws.Range("K33").Value = 159.8
$ws.Range("M33").Value = 69.19999999999999

$ws.Range("H42").Value = 3671.6667
$ws.Range("J42").Value = 6277
$ws.Range("L42").Value = 18831
$ws.Range("N42").Value = -19291

$ws.Range("H92").Value = 1498
$ws.Range("I92").Value = 1497
$ws.Range("K92").Value = 1497
$ws.Range("M92").Value = -249

$ws.Range("H96").Value = 7962.6
$ws.Range("J96").Value = 1878.6
$ws.Range("L96").Value = 5635.799999999999
$ws.Range("N96").Value = -8381.799999999999

$ws.Range("H100").Value = 1000
$ws.Range("J100").Value = 1000
$ws.Range("L100").Value = 1000
$ws.Range("N100").Value = -2082

$ws.Range("H101").Value = 14288433
$ws.Range("J101").Value = 413.5
$ws.Range("L101").Value = 1240.5
$ws.Range("N101").Value = -4484.5

$ws.Range("H113").Value = 2583
$ws.Range("I113").Value = 1743.7778
$ws.Range("K113").Value = 1743.7778
$ws.Range("M113").Value = 1510.2222

$ws.Range("H127").Value = 2180.875
$ws.Range("I127").Value = 2310.2
$ws.Range("K127").Value = 6930.599999999999
$ws.Range("M127").Value = -1970.599999999999

$ws.Range("H132").Value = 2876.7
$ws.Range("I132").Value = 2876.7
$ws.Range("K132").Value = 8630.099999999999
$ws.Range("M132").Value = -6100.099999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5504786.5
$ws.Range("I32").Value = 5389770
$ws.Range("K32").Value = 5389770
$ws.Range("M32").Value = -5389483

$ws.Range("H45").Value = 3756.5
$ws.Range("I45").Value = 3951.6667
$ws.Range("K45").Value = 3951.6667
$ws.Range("M45").Value = -3574.6667

$ws.Range("H61").Value = 2499.8333
$ws.Range("I61").Value = 2549.75
$ws.Range("K61").Value = 2549.75
$ws.Range("M61").Value = -2337.75

$ws.Range("H88").Value = 1641.3636
$ws.Range("I88").Value = 1418.8
$ws.Range("J88").Value = 1826.8334
$ws.Range("K88").Value = 1418.8
$ws.Range("L88").Value = 1826.8334
$ws.Range("M88").Value = -1012.8
$ws.Range("N88").Value = -2638.8334

$ws.Range("H91").Value = 1641.3636
$ws.Range("I91").Value = 1418.8
$ws.Range("J91").Value = 1826.8334
$ws.Range("K91").Value = 1418.8
$ws.Range("L91").Value = 1826.8334
$ws.Range("M91").Value = -14.79999999999995
$ws.Range("N91").Value = -4634.8334

$ws.Range("H102").Value = 1999
$ws.Range("I102").Value = 1999
$ws.Range("K102").Value = 1999
$ws.Range("M102").Value = -377

$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344

$ws.Range("H132").Value = 1306.5
$ws.Range("I132").Value = 1275.3334
$ws.Range("J132").Value = 1400
$ws.Range("K132").Value = 3826.0002
$ws.Range("L132").Value = 4200
$ws.Range("M132").Value = -1296.0002
$ws.Range("N132").Value = -9260

$ws.Range("H136").Value = 2499.8333
$ws.Range("I136").Value = 2549.75
$ws.Range("K136").Value = 7649.25
$ws.Range("M136").Value = -5099.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 19500
$ws.Range("J16").Value = 19500
$ws.Range("L16").Value = 19500
$ws.Range("N16").Value = -19840

$ws.Range("H86").Value = 1091.619
$ws.Range("I86").Value = 890.1539
$ws.Range("K86").Value = 890.1539
$ws.Range("M86").Value = 232.8461

$ws.Range("H89").Value = 1091.619
$ws.Range("I89").Value = 890.1539
$ws.Range("K89").Value = 4450.7695
$ws.Range("M89").Value = 1165.2305

$ws.Range("H99").Value = 1988.3334
$ws.Range("J99").Value = 1988.3334
$ws.Range("L99").Value = 1988.3334
$ws.Range("N99").Value = -4984.3334

$ws.Range("H107").Value = 1707.091
$ws.Range("I107").Value = 1377.9
$ws.Range("K107").Value = 1377.9
$ws.Range("M107").Value = 542.0999999999999

$ws.Range("H134").Value = 2306.8333
$ws.Range("I134").Value = 2306.8333
$ws.Range("K134").Value = 6920.499899999999
$ws.Range("M134").Value = -4385.499899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3235

$ws.Range("H58").Value = 1664.2
$ws.Range("I58").Value = 1330
$ws.Range("J58").Value = 2165.5
$ws.Range("K58").Value = 1330
$ws.Range("L58").Value = 2165.5
$ws.Range("M58").Value = -1127
$ws.Range("N58").Value = -2571.5

$ws.Range("H86").Value = 8412.5
$ws.Range("I86").Value = 7216
$ws.Range("K86").Value = 7216
$ws.Range("M86").Value = -6093

$ws.Range("H89").Value = 8412.5
$ws.Range("I89").Value = 7216
$ws.Range("K89").Value = 36080
$ws.Range("M89").Value = -30464

$ws.Range("H94").Value = 104393
$ws.Range("I94").Value = 223649
$ws.Range("J94").Value = 5013
$ws.Range("K94").Value = 223649
$ws.Range("L94").Value = 5013
$ws.Range("M94").Value = -223198
$ws.Range("N94").Value = -5915

$ws.Range("H113").Value = 3235

$ws.Range("H132").Value = 7139.8
$ws.Range("I132").Value = 6822.1113
$ws.Range("K132").Value = 20466.3339
$ws.Range("M132").Value = -17936.3339

$ws.Range("H134").Value = 2933.3333
$ws.Range("I134").Value = 2900
$ws.Range("K134").Value = 8700
$ws.Range("M134").Value = -6165

$ws.Range("H136").Value = 1664.2
$ws.Range("I136").Value = 1330
$ws.Range("J136").Value = 2165.5
$ws.Range("K136").Value = 3990
$ws.Range("L136").Value = 6496.5
$ws.Range("M136").Value = -1440
$ws.Range("N136").Value = -11596.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 368552
$ws.Range("I2").Value = 366722
$ws.Range("K2").Value = 2200332
$ws.Range("M2").Value = -2200219

$ws.Range("H81").Value = 2182.5
$ws.Range("J81").Value = 2182.5
$ws.Range("L81").Value = 6547.5
$ws.Range("N81").Value = -8793.5

$ws.Range("H84").Value = 2182.5
$ws.Range("J84").Value = 2182.5
$ws.Range("L84").Value = 19642.5
$ws.Range("N84").Value = -30874.5

$ws.Range("H118").Value = 2399
$ws.Range("I118").Value = 2399
$ws.Range("K118").Value = 7197
$ws.Range("M118").Value = -5954

$ws.Range("H128").Value = 499934.5
$ws.Range("I128").Value = 499934.5
$ws.Range("K128").Value = 1499803.5
$ws.Range("M128").Value = -1494823.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2118.6667
$ws.Range("I102").Value = 2105.75
$ws.Range("K102").Value = 2105.75
$ws.Range("M102").Value = -483.75

$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524

$ws.Range("H113").Value = 2833
$ws.Range("J113").Value = 3999.5
$ws.Range("L113").Value = 3999.5
$ws.Range("N113").Value = -8339.5

$ws.Range("H132").Value = 5995.0625
$ws.Range("I132").Value = 5995.0625
$ws.Range("K132").Value = 17985.1875
$ws.Range("M132").Value = -15455.1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1828.1666
$ws.Range("I16").Value = 1828.1666
$ws.Range("K16").Value = 1828.1666
$ws.Range("M16").Value = -1658.1666

$ws.Range("H40").Value = 4659.2
$ws.Range("I40").Value = 4249
$ws.Range("K40").Value = 4249
$ws.Range("M40").Value = -4113

$ws.Range("H100").Value = 3899.3333
$ws.Range("I100").Value = 4068.6924
$ws.Range("J100").Value = 2798.5
$ws.Range("K100").Value = 4068.6924
$ws.Range("L100").Value = 2798.5
$ws.Range("M100").Value = -3527.6924
$ws.Range("N100").Value = -3880.5

$ws.Range("H132").Value = 3405.4614
$ws.Range("I132").Value = 3198.4443
$ws.Range("K132").Value = 9595.332900000001
$ws.Range("M132").Value = -7065.332900000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 16054.429
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 16054.429
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 16054.429
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -17036.429

$ws.Range("H81").Value = 1112227.9
$ws.Range("I81").Value = 1256.125
$ws.Range("J81").Value = 10000002
$ws.Range("K81").Value = 2512.25
$ws.Range("L81").Value = 20000004
$ws.Range("M81").Value = -1451.25
$ws.Range("N81").Value = -20002126

$ws.Range("H84").Value = 1112227.9
$ws.Range("I84").Value = 1256.125
$ws.Range("J84").Value = 10000002
$ws.Range("K84").Value = 12561.25
$ws.Range("L84").Value = 100000020
$ws.Range("M84").Value = -7257.25
$ws.Range("N84").Value = -100010628

$ws.Range("H100").Value = 8334074.5
$ws.Range("I100").Value = 11111806
$ws.Range("K100").Value = 22223612
$ws.Range("M100").Value = -22223071
